$wb = $excel.ActiveWorkbook

# Sheet 1: départements
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(104, 3).Value = 6.976744186046512
$ws1.Cells.Item(104, 4).Value = 9
$ws1.Cells.Item(104, 5).Value = 129
$ws1.Cells.Item(152, 3).Value = 15.78947368421053
$ws1.Cells.Item(152, 4).Value = 6
$ws1.Cells.Item(298, 3).Value = 2.054794520547945
$ws1.Cells.Item(298, 4).Value = 3
$ws1.Cells.Item(362, 3).Value = 3.015075376884422
$ws1.Cells.Item(362, 4).Value = 6
$ws1.Cells.Item(421, 3).Value = 35.59322033898305
$ws1.Cells.Item(421, 4).Value = 21
$ws1.Cells.Item(449, 3).Value = 49.25373134328358
$ws1.Cells.Item(449, 4).Value = 33
$ws1.Cells.Item(454, 3).Value = 46.66666666666666
$ws1.Cells.Item(454, 4).Value = 14
$ws1.Cells.Item(454, 5).Value = 30
$ws1.Cells.Item(492, 3).Value = 22.58064516129032
$ws1.Cells.Item(492, 5).Value = 62
$ws1.Cells.Item(527, 3).Value = 13.95348837209302
$ws1.Cells.Item(527, 4).Value = 6
$ws1.Cells.Item(551, 3).Value = 26.53061224489796
$ws1.Cells.Item(551, 4).Value = 13
$ws1.Cells.Item(551, 5).Value = 49
$ws1.Cells.Item(686, 3).Value = 3.424657534246575
$ws1.Cells.Item(686, 4).Value = 5
$ws1.Cells.Item(750, 3).Value = 6.5
$ws1.Cells.Item(750, 4).Value = 13

# Sheet 2: régions
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(40, 4).Value = 50.91
$ws2.Cells.Item(40, 5).Value = 112
$ws2.Cells.Item(56, 4).Value = 7.32
$ws2.Cells.Item(56, 5).Value = 32
$ws2.Cells.Item(75, 4).Value = 17.78
$ws2.Cells.Item(75, 5).Value = 80
$ws2.Cells.Item(75, 6).Value = 450
$ws2.Cells.Item(76, 4).Value = 40.23
$ws2.Cells.Item(76, 5).Value = 107
$ws2.Cells.Item(76, 6).Value = 266
$ws2.Cells.Item(85, 4).Value = 43.07
$ws2.Cells.Item(85, 5).Value = 115
$ws2.Cells.Item(95, 4).Value = 1.8
$ws2.Cells.Item(95, 5).Value = 24
$ws2.Cells.Item(98, 4).Value = 3.14
$ws2.Cells.Item(98, 5).Value = 42
$ws2.Cells.Item(101, 4).Value = 6.28
$ws2.Cells.Item(101, 5).Value = 44
$ws2.Cells.Item(101, 6).Value = 701
$ws2.Cells.Item(102, 4).Value = 16.1
$ws2.Cells.Item(102, 6).Value = 385
$ws2.Cells.Item(104, 4).Value = 1.99
$ws2.Cells.Item(104, 5).Value = 16
$ws2.Cells.Item(107, 4).Value = 3.34
$ws2.Cells.Item(107, 5).Value = 27

# Sheet 3: national
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(2, 2).Value = 8.55
$ws3.Cells.Item(2, 4).Value = 7715
$ws3.Cells.Item(3, 2).Value = 21
$ws3.Cells.Item(3, 3).Value = 882
$ws3.Cells.Item(4, 2).Value = 45.82
$ws3.Cells.Item(4, 3).Value = 1151
$ws3.Cells.Item(4, 4).Value = 2512
